# Auto-generated edit script applying the diff to Sophia_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 257.14285
$ws.Range("H19").Value = 1507.3334
$ws.Range("J19").Value = 1540
$ws.Range("L19").Value = 1540
$ws.Range("N19").Value = -1890
$ws.Range("H111").Value = 416.66666
$ws.Range("J111").Value = 425
$ws.Range("L111").Value = 1275
$ws.Range("N111").Value = -7409
$ws.Range("H112").Value = 2679.842
$ws.Range("J112").Value = 2693.2222
$ws.Range("L112").Value = 8079.6666
$ws.Range("N112").Value = -10295.6666
$ws.Range("H127").Value = 450
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5435.2856
$ws.Range("I61").Value = 5341.1665
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 5341.1665
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -5129.1665
$ws.Range("N61").Value = -6424
$ws.Range("H74").Value = 5056.091
$ws.Range("I74").Value = 5579.154
$ws.Range("K74").Value = 5579.154
$ws.Range("M74").Value = -4705.154
$ws.Range("H77").Value = 5056.091
$ws.Range("I77").Value = 5579.154
$ws.Range("K77").Value = 27895.77
$ws.Range("M77").Value = -23527.77
$ws.Range("H97").Value = 961.44446
$ws.Range("I97").Value = 950.4286
$ws.Range("K97").Value = 950.4286
$ws.Range("M97").Value = -454.4286
$ws.Range("H102").Value = 2628.5
$ws.Range("I102").Value = 2338.6667
$ws.Range("J102").Value = 3498
$ws.Range("K102").Value = 2338.6667
$ws.Range("L102").Value = 3498
$ws.Range("M102").Value = -716.6667000000002
$ws.Range("N102").Value = -6742
$ws.Range("H136").Value = 5435.2856
$ws.Range("I136").Value = 5341.1665
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 16023.4995
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -13473.4995
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14133
$ws.Range("I20").Value = 2200
$ws.Range("J20").Value = 37999
$ws.Range("K20").Value = 2200
$ws.Range("L20").Value = 37999
$ws.Range("M20").Value = -1953
$ws.Range("N20").Value = -38493
$ws.Range("H94").Value = 10001.6
$ws.Range("I94").Value = 10004
$ws.Range("K94").Value = 10004
$ws.Range("M94").Value = -9553

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 973.75
$ws.Range("I31").Value = 941.4286
$ws.Range("J31").Value = 1200
$ws.Range("K31").Value = 941.4286
$ws.Range("L31").Value = 1200
$ws.Range("M31").Value = -646.4286
$ws.Range("N31").Value = -1790
$ws.Range("H34").Value = 973.75
$ws.Range("I34").Value = 941.4286
$ws.Range("J34").Value = 1200
$ws.Range("K34").Value = 941.4286
$ws.Range("L34").Value = 1200
$ws.Range("M34").Value = -739.4286
$ws.Range("N34").Value = -1604

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 178.6
$ws.Range("I23").Value = 79
$ws.Range("J23").Value = 203.5
$ws.Range("K23").Value = 237
$ws.Range("L23").Value = 610.5
$ws.Range("M23").Value = -2
$ws.Range("N23").Value = -1080.5
$ws.Range("H104").Value = 5185.5713
$ws.Range("I104").Value = 4000
$ws.Range("J104").Value = 5244.85
$ws.Range("K104").Value = 12000
$ws.Range("L104").Value = 15734.55
$ws.Range("M104").Value = -9379
$ws.Range("N104").Value = -20976.55
$ws.Range("H109").Value = 699
$ws.Range("I109").Value = 699
$ws.Range("K109").Value = 2097
$ws.Range("M109").Value = -1057
$ws.Range("H113").Value = 1371.4286
$ws.Range("J113").Value = 1760
$ws.Range("L113").Value = 5280
$ws.Range("N113").Value = -9620
$ws.Range("H121").Value = 1431.5
$ws.Range("J121").Value = 4000
$ws.Range("L121").Value = 12000
$ws.Range("N121").Value = -14620
$ws.Range("H134").Value = 4766.6665
$ws.Range("I134").Value = 4766.6665
$ws.Range("K134").Value = 14299.9995
$ws.Range("M134").Value = -9229.999500000002
$ws.Range("H137").Value = 875
$ws.Range("I137").Value = 750
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 2250
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = 2850
$ws.Range("N137").Value = -13200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 6000000
$ws.Range("J10").Value = 6000000
$ws.Range("L10").Value = 6000000
$ws.Range("N10").Value = -6000338
$ws.Range("H13").Value = 2376.25
$ws.Range("I13").Value = 2168.3333
$ws.Range("K13").Value = 2168.3333
$ws.Range("M13").Value = -2029.3333
$ws.Range("H70").Value = 3999
$ws.Range("J70").Value = 3999
$ws.Range("L70").Value = 3999
$ws.Range("N70").Value = -4539
$ws.Range("H73").Value = 3999
$ws.Range("J73").Value = 3999
$ws.Range("L73").Value = 3999
$ws.Range("N73").Value = -5871
$ws.Range("H122").Value = 3651.5
$ws.Range("I122").Value = 3602
$ws.Range("J122").Value = 3800
$ws.Range("K122").Value = 10806
$ws.Range("L122").Value = 11400
$ws.Range("M122").Value = -8356
$ws.Range("N122").Value = -16300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1750
$ws.Range("I7").Value = 1750
$ws.Range("K7").Value = 1750
$ws.Range("M7").Value = -1638
$ws.Range("H22").Value = 9627
$ws.Range("I22").Value = 14287
$ws.Range("J22").Value = 7508.8184
$ws.Range("K22").Value = 14287
$ws.Range("L22").Value = 7508.8184
$ws.Range("M22").Value = -13992
$ws.Range("N22").Value = -8098.8184
$ws.Range("H27").Value = 9627
$ws.Range("I27").Value = 14287
$ws.Range("J27").Value = 7508.8184
$ws.Range("K27").Value = 14287
$ws.Range("L27").Value = 7508.8184
$ws.Range("M27").Value = -14180
$ws.Range("N27").Value = -7722.8184
$ws.Range("H122").Value = 4601.3335
$ws.Range("I122").Value = 3904
$ws.Range("K122").Value = 11712
$ws.Range("M122").Value = -9262
$ws.Range("H126").Value = 1750
$ws.Range("I126").Value = 1750
$ws.Range("K126").Value = 5250
$ws.Range("M126").Value = -2780
$ws.Range("H132").Value = 5944.077
$ws.Range("I132").Value = 3152
$ws.Range("K132").Value = 9456
$ws.Range("M132").Value = -6926

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 16111
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 16111
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 16111
$ws.Range("N21").Value = -16581
$ws.Range("H29").Value = 10000
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10580
$ws.Range("H35").Value = 16111
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 16111
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 16111
$ws.Range("N35").Value = -16691
$ws.Range("H81").Value = 3862.4707
$ws.Range("I81").Value = 1721.8334
$ws.Range("K81").Value = 3443.6668
$ws.Range("M81").Value = -2382.6668
$ws.Range("H84").Value = 3862.4707
$ws.Range("I84").Value = 1721.8334
$ws.Range("K84").Value = 17218.334
$ws.Range("M84").Value = -11914.334
$ws.Range("M21").ClearContents()
$ws.Range("M35").ClearContents()
